$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (column D) values
$ws.Range("D2").Value = "67.823.26"
$ws.Range("D3").Value = "3.321.96"
$ws.Range("D5").Value = "'578.84"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'174.98"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'0.588"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "3.319.15"
$ws.Range("D11").Value = "'0.575"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "'45.37"
$ws.Range("D12").Style = "Normal"
$ws.Range("D14").Value = "'661.88"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "3.864.09"
$ws.Range("D16").Value = "'8.40"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "67.737.86"
$ws.Range("D19").Value = "3.324.00"
$ws.Range("D24").Value = "'16.87"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Value = "'98.15"
$ws.Range("D25").Style = "Normal"
$ws.Range("D28").Value = "'9.23"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Value = "'33.23"
$ws.Range("D29").Style = "Normal"
$ws.Range("D31").Value = "'7.24"
$ws.Range("D31").Style = "Normal"
$ws.Range("D35").Value = "'0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Value = "3.679.12"
$ws.Range("D38").Value = "'3.22"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Value = "'34.37"
$ws.Range("D39").Style = "Normal"
$ws.Range("D42").Value = "'3.10"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Value = "'3.31"
$ws.Range("D43").Style = "Normal"
$ws.Range("D47").Value = "'2.58"
$ws.Range("D47").Style = "Normal"
$ws.Range("D51").Value = "'127.66"
$ws.Range("D51").Style = "Normal"

# Update Volume(1h) (column E) values
$ws.Range("E2").Value = "  -0.59%  "
$ws.Range("E3").Value = "  -0.01%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("E5").Value = "  -0.77%  "
$ws.Range("E6").Value = "  -4.45%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  -0.69%  "
$ws.Range("E9").Value = "  +0.01%  "
$ws.Range("E10").Value = "  -0.72%  "
$ws.Range("E11").Value = "  -0.95%  "
$ws.Range("E12").Value = "  -2.10%  "
$ws.Range("E13").Value = "  -2.63%  "
$ws.Range("E14").Value = "  +3.89%  "
$ws.Range("E15").Value = "  +0.25%  "
$ws.Range("E16").Value = "  -0.85%  "
$ws.Range("E17").Value = "  -0.80%  "
$ws.Range("E18").Value = "  -0.90%  "
$ws.Range("E19").Value = "  +0.06%  "
$ws.Range("E20").Value = "  -1.88%  "
$ws.Range("E21").Value = "  -0.05%  "
$ws.Range("E22").Value = "  -1.94%  "
$ws.Range("E23").Value = "  +5.05%  "
$ws.Range("E24").Value = "  -4.61%  "
$ws.Range("E25").Value = "  +0.98%  "
$ws.Range("E26").Value = "  -3.98%  "
$ws.Range("E27").Value = "  -4.21%  "
$ws.Range("E28").Value = "  -3.82%  "
$ws.Range("E29").Value = "  +2.38%  "
$ws.Range("E30").Value = "  -2.38%  "
$ws.Range("E31").Value = "  +8.06%  "
$ws.Range("E32").Value = "  -6.20%  "
$ws.Range("E33").Value = "  -0.51%  "
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("E35").Value = "  -0.08%  "
$ws.Range("E36").Value = "  -6.95%  "
$ws.Range("E37").Value = "  +0.69%  "
$ws.Range("E38").Value = "  -8.29%  "
$ws.Range("E39").Value = "  +5.28%  "
$ws.Range("E40").Value = "  +0.63%  "
$ws.Range("E41").Value = "  -2.54%  "
$ws.Range("E42").Value = "  -5.09%  "
$ws.Range("E43").Value = "  -2.46%  "
$ws.Range("E44").Value = "  -1.59%  "
$ws.Range("E45").Value = "  -4.10%  "
$ws.Range("E46").Value = "  -2.60%  "
$ws.Range("E47").Value = "  +1.19%  "
$ws.Range("E48").Value = "  -0.99%  "
$ws.Range("E49").Value = "  -0.33%  "
$ws.Range("E50").Value = "  -3.60%  "
$ws.Range("E51").Value = "  -2.80%  "
